# Auto-generated edit script applying scheduled-runner price/profit refresh
# to the Leve-crafting-profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H55").Value = 153.83333
$ws.Range("I55").Value = 101.1
$ws.Range("J55").Value = 219.75
$ws.Range("K55").Value = 101.1
$ws.Range("L55").Value = 219.75
$ws.Range("M55").Value = 112.9
$ws.Range("N55").Value = -647.75

$ws.Range("H92").Value = 841.1818
$ws.Range("I92").Value = 875.3
$ws.Range("K92").Value = 875.3
$ws.Range("M92").Value = 372.7

$ws.Range("H137").Value = 685.2857
$ws.Range("I137").Value = 699.6667
$ws.Range("J137").Value = 599
$ws.Range("K137").Value = 2099.0001
$ws.Range("L137").Value = 1797
$ws.Range("M137").Value = 450.9998999999998
$ws.Range("N137").Value = -6897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 600
$ws.Range("I15").Value = 600
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -250
$ws.Range("N15").ClearContents()

$ws.Range("H32").Value = 2730.2083
$ws.Range("I32").Value = 2568.0908
$ws.Range("J32").Value = 4513.5
$ws.Range("K32").Value = 2568.0908
$ws.Range("L32").Value = 4513.5
$ws.Range("M32").Value = -2281.0908
$ws.Range("N32").Value = -5087.5

$ws.Range("H74").Value = 2084.375
$ws.Range("I74").Value = 2023.3334
$ws.Range("K74").Value = 2023.3334
$ws.Range("M74").Value = -1149.3334

$ws.Range("H77").Value = 2084.375
$ws.Range("I77").Value = 2023.3334
$ws.Range("K77").Value = 10116.667
$ws.Range("M77").Value = -5748.666999999999

$ws.Range("H122").Value = 1510.75
$ws.Range("I122").Value = 1510.75
$ws.Range("K122").Value = 4532.25
$ws.Range("M122").Value = -2082.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2136.7693
$ws.Range("I99").Value = 1811
$ws.Range("J99").Value = 2869.75
$ws.Range("K99").Value = 1811
$ws.Range("L99").Value = 2869.75
$ws.Range("M99").Value = -313
$ws.Range("N99").Value = -5865.75

$ws.Range("H134").Value = 1172.5
$ws.Range("I134").Value = 1172.5
$ws.Range("K134").Value = 3517.5
$ws.Range("M134").Value = -982.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1389.4762
$ws.Range("I7").Value = 898.26666
$ws.Range("J7").Value = 2617.5
$ws.Range("K7").Value = 898.26666
$ws.Range("L7").Value = 2617.5
$ws.Range("M7").Value = -785.26666
$ws.Range("N7").Value = -2843.5

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H86").Value = 749751.5
$ws.Range("I86").Value = 998002
$ws.Range("K86").Value = 998002
$ws.Range("M86").Value = -996879

$ws.Range("H89").Value = 749751.5
$ws.Range("I89").Value = 998002
$ws.Range("K89").Value = 4990010
$ws.Range("M89").Value = -4984394

$ws.Range("H94").Value = 2559.8333
$ws.Range("I94").Value = 2478.6667
$ws.Range("K94").Value = 2478.6667
$ws.Range("M94").Value = -2027.6667

$ws.Range("H99").Value = 2572.2222
$ws.Range("I99").Value = 2572.2222
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2572.2222
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1074.2222
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 5046.75
$ws.Range("I105").Value = 2074.8
$ws.Range("K105").Value = 2074.8
$ws.Range("M105").Value = -327.8000000000002

$ws.Range("H122").Value = 1131.75
$ws.Range("I122").Value = 1131.75
$ws.Range("K122").Value = 3395.25
$ws.Range("M122").Value = -945.25

$ws.Range("H126").Value = 2572.2222
$ws.Range("I126").Value = 2572.2222
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7716.6666
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5246.6666
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 278.2857
$ws.Range("I33").Value = 216
$ws.Range("K33").Value = 1296
$ws.Range("M33").Value = -1013

$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -754

$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -49868

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H128").Value = 277417.5
$ws.Range("I128").Value = 277417.5
$ws.Range("K128").Value = 832252.5
$ws.Range("M128").Value = -827272.5

$ws.Range("H131").Value = 1053.3334
$ws.Range("J131").Value = 1053.3334
$ws.Range("L131").Value = 3160.0002
$ws.Range("N131").Value = -13240.0002

$ws.Range("H137").Value = 3919
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3919
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 11757
$ws.Range("N137").Value = -21957
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4496.3335
$ws.Range("I97").Value = 4494
$ws.Range("J97").Value = 4497.5
$ws.Range("K97").Value = 4494
$ws.Range("L97").Value = 4497.5
$ws.Range("M97").Value = -3998
$ws.Range("N97").Value = -5489.5

$ws.Range("H132").Value = 3999
$ws.Range("I132").Value = 3999
$ws.Range("K132").Value = 11997
$ws.Range("M132").Value = -9467

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4500
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H40").Value = 4244.75
$ws.Range("I40").Value = 3993
$ws.Range("K40").Value = 3993
$ws.Range("M40").Value = -3857

$ws.Range("H82").Value = 1192.5625
$ws.Range("I82").Value = 1198.9231
$ws.Range("J82").Value = 1165
$ws.Range("K82").Value = 1198.9231
$ws.Range("L82").Value = 1165
$ws.Range("M82").Value = -837.9231
$ws.Range("N82").Value = -1887

$ws.Range("H85").Value = 1192.5625
$ws.Range("I85").Value = 1198.9231
$ws.Range("J85").Value = 1165
$ws.Range("K85").Value = 1198.9231
$ws.Range("L85").Value = 1165
$ws.Range("M85").Value = 49.07690000000002
$ws.Range("N85").Value = -3661

$ws.Range("H126").Value = 4500
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 6579
$ws.Range("I132").Value = 7666
$ws.Range("J132").Value = 4948.5
$ws.Range("K132").Value = 22998
$ws.Range("L132").Value = 14845.5
$ws.Range("M132").Value = -20468
$ws.Range("N132").Value = -19905.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 39998
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 39998
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 39998
$ws.Range("N47").Value = -41142
$ws.Range("M47").ClearContents()

$ws.Range("H48").Value = 44000
$ws.Range("J48").Value = 44000
$ws.Range("L48").Value = 44000
$ws.Range("N48").Value = -45138

$ws.Range("H126").Value = 2803
$ws.Range("I126").Value = 2803
$ws.Range("K126").Value = 8409
$ws.Range("M126").Value = -5939
